$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename columns, add "Year" header in D1 ---
$ws.Range("A1").Value = "Rank"
$ws.Range("B1").Value = "City Name"
$ws.Range("C1").Value = "Overnight International Visitor Spend (US`$ bn)"
$ws.Range("C1").Style = "Normal"
$ws.Range("D1").Value = "Year"

# --- Column D: fill rows 2-21 with the year 2011 (replacing the blank placeholder cells) ---
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 4).Value = 2011
}

# --- Remove the now-unused trailing rows 22-24 (only held placeholder text before) ---
$ws.Rows("22:24").Delete()

# --- Column C width: widen to fit the new, longer header text ---
$ws.Columns.Item(3).ColumnWidth = 35.83

# --- Selection matches the authored file: D2:D21 with active cell D2 ---
$null = $ws.Range("D2:D21").Select()
